$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the SVM prediction results (netfit) for the Training row
$ws.Range("E3").Value = 0.98809519999999995
$ws.Range("F3").Value = 0.91666669999999995
$ws.Range("G3").Value = 0.96031750000000005

# Update the active cell selection to match the saved state
$ws.Range("G10").Select()
